# Change font size for buttons (the "L" and "R" shapes) to 32pt,
# matching the numbered circle buttons (1-5) already at 32pt.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$g = $s.Shapes.Item(1)

for ($i = 1; $i -le $g.GroupItems.Count; $i++) {
    $sh = $g.GroupItems.Item($i)
    if ($sh.HasTextFrame) {
        $txt = $sh.TextFrame.TextRange.Text
        if ($txt -eq "L" -or $txt -eq "R") {
            $sh.TextFrame.TextRange.Font.Size = 32
        }
    }
}
